# Vessel-Fleet-Optimization "Inputs" workbook update
# ---------------------------------------------------
# The spare_parts master list is trimmed from 4 parts (S1..S4) down to 2
# parts (S1, S2). All the dependent sheets that spill a
# FILTER(spare_parts!A2:A100, ...) formula down column A need their spill
# shrunk from A2:A5 to A2:A3, and a handful of the manually entered
# quantities are updated as part of the same pass (better test data for
# the plotter / verification run).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) spare_parts: delete the S3 / S4 rows outright
# ---------------------------------------------------------------------
$wsSpareParts = $wb.Worksheets.Item("spare_parts")
$wsSpareParts.Range("A4:C5").ClearContents()

# ---------------------------------------------------------------------
# 2) holding_costs: shrink the spilled FILTER, drop the now-empty S3/S4
#    label cells, and update the entered cost values
# ---------------------------------------------------------------------
$wsHoldingCosts = $wb.Worksheets.Item("holding_costs")
$wsHoldingCosts.Range("A2:A3").FormulaArray = "=_xlfn._xlws.FILTER(spare_parts!A2:A100, spare_parts!A2:A100<>"""")"
$wsHoldingCosts.Range("A4:A5").ClearContents()
$wsHoldingCosts.Range("B2").Value = 8
$wsHoldingCosts.Range("C2").Value = 5
$wsHoldingCosts.Range("B3").Value = 5
$wsHoldingCosts.Range("C3").Value = 3

# ---------------------------------------------------------------------
# 3) spare_parts_required: shrink the spilled FILTER, drop the S3/S4
#    label cells, and update the required quantities
# ---------------------------------------------------------------------
$wsSparePartsRequired = $wb.Worksheets.Item("spare_parts_required")
$wsSparePartsRequired.Range("A2:A3").FormulaArray = "=_xlfn._xlws.FILTER(spare_parts!A2:A100, spare_parts!A2:A100<>"""")"
$wsSparePartsRequired.Range("A4:A5").ClearContents()
$wsSparePartsRequired.Range("C2").Value = 1
$wsSparePartsRequired.Range("C3").Value = 2

# ---------------------------------------------------------------------
# 4) max_capacity: shrink the spilled FILTER, drop the S3/S4 label
#    cells, and update the capacity values
# ---------------------------------------------------------------------
$wsMaxCapacity = $wb.Worksheets.Item("max_capacity")
$wsMaxCapacity.Range("A2:A3").FormulaArray = "=_xlfn._xlws.FILTER(spare_parts!A2:A100, spare_parts!A2:A100<>"""")"
$wsMaxCapacity.Range("A4:A5").ClearContents()
$wsMaxCapacity.Range("B2").Value = 50
$wsMaxCapacity.Range("C2").Value = 40
$wsMaxCapacity.Range("B3").Value = 60
$wsMaxCapacity.Range("C3").Value = 40

# ---------------------------------------------------------------------
# 5) reorder_level: shrink the spilled FILTER and drop the S3/S4 label
#    cells (the quantities themselves are unchanged)
# ---------------------------------------------------------------------
$wsReorderLevel = $wb.Worksheets.Item("reorder_level")
$wsReorderLevel.Range("A2:A3").FormulaArray = "=_xlfn._xlws.FILTER(spare_parts!A2:A100, spare_parts!A2:A100<>"""")"
$wsReorderLevel.Range("A4:A5").ClearContents()

# ---------------------------------------------------------------------
# 6) Selection / active-sheet bookkeeping to match the saved view state
# ---------------------------------------------------------------------
$wsSpareParts.Range("A4:C5").Select()
$wsHoldingCosts.Range("C3").Select()
$wsSparePartsRequired.Range("C2").Select()
$wsMaxCapacity.Range("I11").Select()

# reorder_level ends up as the active tab/selected cell when the file is
# saved
$wsReorderLevel.Range("B2").Select()
